$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Category renamed from "Cake" to "Sweet Cake" for all three product rows.
$ws.Range("B2").Value = "Sweet Cake"
$ws.Range("B3").Value = "Sweet Cake"
$ws.Range("B4").Value = "Sweet Cake"

# Row heights grew (re-wrapped header/description rows).
$ws.Rows.Item(1).RowHeight = 37.5
$ws.Rows.Item(2).RowHeight = 75
$ws.Rows.Item(3).RowHeight = 75
$ws.Rows.Item(4).RowHeight = 75

# Active selection moved to B4.
[void]$ws.Range("B4").Select()
